$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.408218
$ws.Range("H2").Value = 61.224654
$ws.Range("I2").Value = 0.1108535210972707
$ws.Range("J2").Value = 0.1108535210972707
$ws.Range("M2").Value = 3.135398666666667
$ws.Range("N2").Value = 9.406196000000001
$ws.Range("O2").Value = 0.1723049126704688
$ws.Range("P2").Value = 0.1723049126704688
$ws.Range("Q2").Value = 63.98789950624268
$ws.Range("R2").Value = 575.891095556184
$ws.Range("S2").Value = 0.0191006062718792
$ws.Range("T2").Value = 0.0191006062718792
$ws.Range("G3").Value = 20.408218
$ws.Range("H3").Value = 61.224654
$ws.Range("I3").Value = 0.1108535210972707
$ws.Range("J3").Value = 0.1108535210972707
$ws.Range("O3").Value = 0.1733096678828815
$ws.Range("P3").Value = 0.1733096678828815
$ws.Range("Q3").Value = 64.361029758676
$ws.Range("R3").Value = 579.2492678280839
$ws.Range("S3").Value = 0.01921198692501598
$ws.Range("T3").Value = 0.01921198692501598
$ws.Range("G4").Value = 20.408218
$ws.Range("H4").Value = 61.224654
$ws.Range("I4").Value = 0.1108535210972707
$ws.Range("J4").Value = 0.1108535210972707
$ws.Range("M4").Value = 0.4900660000000001
$ws.Range("N4").Value = 1.470198
$ws.Range("O4").Value = 0.02693143306797965
$ws.Range("P4").Value = 0.02693143306797965
$ws.Range("Q4").Value = 10.001373762388
$ws.Range("R4").Value = 90.01236386149201
$ws.Range("S4").Value = 0.002985444183781016
$ws.Range("T4").Value = 0.002985444183781016
$ws.Range("G5").Value = 20.408218
$ws.Range("H5").Value = 61.224654
$ws.Range("I5").Value = 0.1108535210972707
$ws.Range("J5").Value = 0.1108535210972707
$ws.Range("M5").Value = 11.417657
$ws.Range("N5").Value = 34.252971
$ws.Range("O5").Value = 0.62745398637867
$ws.Range("P5").Value = 0.6274539863786701
$ws.Range("Q5").Value = 233.014033105226
$ws.Range("R5").Value = 2097.126297947034
$ws.Range("S5").Value = 0.06955548371659449
$ws.Range("T5").Value = 0.06955548371659451
$ws.Range("I6").Value = 0.2566851044076959
$ws.Range("J6").Value = 0.256685104407696
$ws.Range("M6").Value = 3.135398666666667
$ws.Range("N6").Value = 9.406196000000001
$ws.Range("O6").Value = 0.1723049126704688
$ws.Range("P6").Value = 0.1723049126704688
$ws.Range("Q6").Value = 148.1661611017013
$ws.Range("R6").Value = 1333.495449915312
$ws.Range("S6").Value = 0.04422810449877822
$ws.Range("T6").Value = 0.04422810449877823
$ws.Range("I7").Value = 0.2566851044076959
$ws.Range("J7").Value = 0.256685104407696
$ws.Range("O7").Value = 0.1733096678828815
$ws.Range("P7").Value = 0.1733096678828815
$ws.Range("S7").Value = 0.04448601019538054
$ws.Range("T7").Value = 0.04448601019538054
$ws.Range("I8").Value = 0.2566851044076959
$ws.Range("J8").Value = 0.256685104407696
$ws.Range("M8").Value = 0.4900660000000001
$ws.Range("N8").Value = 1.470198
$ws.Range("O8").Value = 0.02693143306797965
$ws.Range("P8").Value = 0.02693143306797965
$ws.Range("Q8").Value = 23.158521650984
$ws.Range("R8").Value = 208.426694858856
$ws.Range("S8").Value = 0.006912897708903232
$ws.Range("T8").Value = 0.006912897708903233
$ws.Range("I9").Value = 0.2566851044076959
$ws.Range("J9").Value = 0.256685104407696
$ws.Range("M9").Value = 11.417657
$ws.Range("N9").Value = 34.252971
$ws.Range("O9").Value = 0.62745398637867
$ws.Range("P9").Value = 0.6274539863786701
$ws.Range("Q9").Value = 539.551931450068
$ws.Range("R9").Value = 4855.967383050612
$ws.Range("S9").Value = 0.1610580920046339
$ws.Range("T9").Value = 0.161058092004634
$ws.Range("G10").Value = 85.307233
$ws.Range("H10").Value = 255.921699
$ws.Range("I10").Value = 0.4633725077375833
$ws.Range("J10").Value = 0.4633725077375833
$ws.Range("M10").Value = 3.135398666666667
$ws.Range("N10").Value = 9.406196000000001
$ws.Range("O10").Value = 0.1723049126704688
$ws.Range("P10").Value = 0.1723049126704688
$ws.Range("Q10").Value = 267.4721846052227
$ws.Range("R10").Value = 2407.249661447004
$ws.Range("S10").Value = 0.07984135947962043
$ws.Range("T10").Value = 0.07984135947962043
$ws.Range("G11").Value = 85.307233
$ws.Range("H11").Value = 255.921699
$ws.Range("I11").Value = 0.4633725077375833
$ws.Range("J11").Value = 0.4633725077375833
$ws.Range("O11").Value = 0.1733096678828815
$ws.Range("P11").Value = 0.1733096678828815
$ws.Range("Q11").Value = 269.031885181906
$ws.Range("R11").Value = 2421.286966637154
$ws.Range("S11").Value = 0.08030693542205848
$ws.Range("T11").Value = 0.08030693542205848
$ws.Range("G12").Value = 85.307233
$ws.Range("H12").Value = 255.921699
$ws.Range("I12").Value = 0.4633725077375833
$ws.Range("J12").Value = 0.4633725077375833
$ws.Range("M12").Value = 0.4900660000000001
$ws.Range("N12").Value = 1.470198
$ws.Range("O12").Value = 0.02693143306797965
$ws.Range("P12").Value = 0.02693143306797965
$ws.Range("Q12").Value = 41.80617444737801
$ws.Range("R12").Value = 376.255570026402
$ws.Range("S12").Value = 0.01247928567767661
$ws.Range("T12").Value = 0.01247928567767661
$ws.Range("G13").Value = 85.307233
$ws.Range("H13").Value = 255.921699
$ws.Range("I13").Value = 0.4633725077375833
$ws.Range("J13").Value = 0.4633725077375833
$ws.Range("M13").Value = 11.417657
$ws.Range("N13").Value = 34.252971
$ws.Range("O13").Value = 0.62745398637867
$ws.Range("P13").Value = 0.6274539863786701
$ws.Range("Q13").Value = 974.0087260130809
$ws.Range("R13").Value = 8766.078534117729
$ws.Range("S13").Value = 0.2907449271582278
$ws.Range("T13").Value = 0.2907449271582278
$ws.Range("G14").Value = 31.12938966666666
$ws.Range("H14").Value = 93.38816899999999
$ws.Range("I14").Value = 0.16908886675745
$ws.Range("J14").Value = 0.16908886675745
$ws.Range("M14").Value = 3.135398666666667
$ws.Range("N14").Value = 9.406196000000001
$ws.Range("O14").Value = 0.1723049126704688
$ws.Range("P14").Value = 0.1723049126704688
$ws.Range("Q14").Value = 97.60304685501379
$ws.Range("R14").Value = 878.4274216951241
$ws.Range("S14").Value = 0.02913484242019097
$ws.Range("T14").Value = 0.02913484242019097
$ws.Range("G15").Value = 31.12938966666666
$ws.Range("H15").Value = 93.38816899999999
$ws.Range("I15").Value = 0.16908886675745
$ws.Range("J15").Value = 0.16908886675745
$ws.Range("O15").Value = 0.1733096678828815
$ws.Range("P15").Value = 0.1733096678828815
$ws.Range("Q15").Value = 98.17219586275266
$ws.Range("R15").Value = 883.5497627647738
$ws.Range("S15").Value = 0.02930473534042646
$ws.Range("T15").Value = 0.02930473534042646
$ws.Range("G16").Value = 31.12938966666666
$ws.Range("H16").Value = 93.38816899999999
$ws.Range("I16").Value = 0.16908886675745
$ws.Range("J16").Value = 0.16908886675745
$ws.Range("M16").Value = 0.4900660000000001
$ws.Range("N16").Value = 1.470198
$ws.Range("O16").Value = 0.02693143306797965
$ws.Range("P16").Value = 0.02693143306797965
$ws.Range("Q16").Value = 15.25545547638467
$ws.Range("R16").Value = 137.299099287462
$ws.Range("S16").Value = 0.004553805497618794
$ws.Range("T16").Value = 0.004553805497618794
$ws.Range("G17").Value = 31.12938966666666
$ws.Range("H17").Value = 93.38816899999999
$ws.Range("I17").Value = 0.16908886675745
$ws.Range("J17").Value = 0.16908886675745
$ws.Range("M17").Value = 11.417657
$ws.Range("N17").Value = 34.252971
$ws.Range("O17").Value = 0.62745398637867
$ws.Range("P17").Value = 0.6274539863786701
$ws.Range("Q17").Value = 355.4246938333443
$ws.Range("R17").Value = 3198.822244500099
$ws.Range("S17").Value = 0.1060954834992138
$ws.Range("T17").Value = 0.1060954834992138
